$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" row: the value is now reported explicitly as the text "false"
# (it used to be left blank). A leading apostrophe forces Excel to store the
# literal text "false" instead of auto-converting it to a Boolean; copying the
# formatting back from the neighbouring "Experimental" label cell keeps B7 on
# its original (unprefixed) cell style.
$ws.Range("B7").Value = "'false"
$ws.Range("A7").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# "Date" row: regenerated at a later timestamp.
$ws.Range("B8").Value = "2023-10-09T22:41:16+02:00"
